# Applies the diff described for Corridas/20251025_150353_analisis_sp500/resultados.xlsx
# Sheet "Metricas" (sheet1): relabel / recompute metrics, add a new
#   "Rendimiento_Esperado_Porcentual" row, shifting N_Acciones and
#   Peso_Total_Acciones down by one row (dimension grows to A1:B10).
# Sheet "Acciones_Seleccionadas" (sheet2): add a new "Rendimiento_Porcentual"
#   column (dimension grows to A1:H11) and replace the whole stock table
#   with the new selection / metrics.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Metricas
# ---------------------------------------------------------------------
$wsMetricas = $wb.Worksheets.Item("Metricas")

$metricLabels = @(
    "Rendimiento_Esperado_Log",
    "Rendimiento_Esperado_Porcentual",
    "Volatilidad",
    "Sharpe_Ratio",
    "Riesgo_Varianza",
    "CVaR_Portafolio",
    "Funcion_Objetivo",
    "N_Acciones",
    "Peso_Total_Acciones"
)

$metricValues = @(
    0.01573400932611407,
    0.01585844059520181,
    0.03917389605175114,
    0.3165545744134992,
    0.001534594131873403,
    -0.03382977644364061,
    0.05679026405857873,
    10,
    0.9999999999999998
)

for ($i = 0; $i -lt $metricLabels.Length; $i++) {
    $row = $i + 2
    $wsMetricas.Cells.Item($row, 1).Value = $metricLabels[$i]
    $wsMetricas.Cells.Item($row, 2).Value = $metricValues[$i]
}

# ---------------------------------------------------------------------
# Sheet 2: Acciones_Seleccionadas
# ---------------------------------------------------------------------
$wsAcciones = $wb.Worksheets.Item("Acciones_Seleccionadas")

# New header row (column D "Rendimiento_Porcentual" is inserted, shifting
# the remaining headers one column to the right).
$wsAcciones.Cells.Item(1, 1).Value = "Ticker"
$wsAcciones.Cells.Item(1, 2).Value = "Peso_W"
$wsAcciones.Cells.Item(1, 3).Value = "Rendimiento_Log"
$wsAcciones.Cells.Item(1, 4).Value = "Rendimiento_Porcentual"
$wsAcciones.Cells.Item(1, 5).Value = "Desvio_Estandar"
$wsAcciones.Cells.Item(1, 6).Value = "VaR_95"
$wsAcciones.Cells.Item(1, 7).Value = "CVaR_95"
$wsAcciones.Cells.Item(1, 8).Value = "Prob_Perdida"

# New stock selection data (Ticker, Peso_W, Rendimiento_Log,
# Rendimiento_Porcentual, Desvio_Estandar, VaR_95, CVaR_95, Prob_Perdida)
$accionesData = @(
    @("RSG", 0.3,                0.01593411162495392, 0.01606173654377585, 0.04897055828025669, -0.06009945731390331, -0.08589280382996545, 0.3652062211561631),
    @("MCD", 0.2485366752572068, 0.01105134070186723, 0.01111263234437421, 0.0511631372545762,  -0.05953797337572052, -0.08479405149596911, 0.378843061609838),
    @("MMC", 0.1014633247427929, 0.01115662774737973, 0.01121909501058549, 0.05354270517980156, -0.06283680939995701, -0.09090181999418419, 0.3492962149470633),
    @("AJG", 0.05,               0.01577129163518248, 0.01589631484956922, 0.05840188965830416, -0.06845050899881583, -0.09831927763281612, 0.3640652192290516),
    @("COST", 0.05,              0.02031154968988277, 0.02051923295500946, 0.06575729349546175, -0.06992220936854815, -0.09975833595823232, 0.3523794465681823),
    @("CTAS", 0.05,              0.02303776214964436, 0.02330518102099788, 0.06352659485413112, -0.07606359437629205, -0.1095501084668015,  0.3424831250497379),
    @("LLY", 0.05,               0.03532633271988238, 0.03595772054510604, 0.07642567826879772, -0.08116955916157426, -0.116061859619291,   0.3653204385164766),
    @("PGR", 0.05,               0.01856688928963813, 0.01874032570732598, 0.05532539794482709, -0.07526015336134566, -0.1078942195780875,  0.3652459514316563),
    @("TJX", 0.05,               0.01631437673927512, 0.01644818284754712, 0.05931325343532361, -0.07222838272756965, -0.1033520595563616,  0.3634465622327827),
    @("WRB", 0.05,               0.01217427416218126, 0.01224868228624798, 0.05513970027329756, -0.06722891743011217, -0.09645117124281785, 0.3777689387914418)
)

for ($i = 0; $i -lt $accionesData.Length; $i++) {
    $row = $i + 2
    $values = $accionesData[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 1
        $wsAcciones.Cells.Item($row, $col).Value = $values[$j]
    }
}
